# Chambers/White material_and_methods workbook:
# - study_id (column A) values for rows 4-6 are renamed to reflect the
#   split of "Chambers_et_al_2020_inland/bay" into a single
#   "Chambers_et_al_2019" study, and "White_et_al_2019" into the two
#   "White_et_al_2020_a" / "White_et_al_2020_b" studies.
# - the active sheet's scroll position / selection moved down a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("material_and_methods")

# Assign A6 first, then A4, then A5 so the workbook's shared-string table
# is rebuilt in the same order as the target file
# (Chambers_et_al_2019, White_et_al_2020_a, White_et_al_2020_b).
$ws.Range("A6").Value = "Chambers_et_al_2019"
$ws.Range("A4").Value = "White_et_al_2020_a"
$ws.Range("A5").Value = "White_et_al_2020_b"

# Scroll the sheet down one row and move the selection to A5, matching the
# saved view state (topLeftCell A2 / active cell A5).
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("A5").Select()
